# Apply INE data update: rename municipality "Alcadozo" -> "Abengibre",
# set proper Genero labels (HOMBRE/MUJER) for the per-gender rows, and
# update the Valor figures to match the new municipality's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: Municipio - rows 2 through 13 change from Alcadozo to Abengibre
$ws.Range("A2:A13").Value = "Abengibre"

# Column B: Genero labels for rows 6-9 (HOMBRE) and 10-13 (MUJER)
$ws.Range("B6:B9").Value = "HOMBRE"
$ws.Range("B10:B13").Value = "MUJER"

# Column D: Valor - updated figures
$ws.Range("D2").Value = 760
$ws.Range("D3").Value = 739
$ws.Range("D4").Value = 748
$ws.Range("D5").Value = 761
$ws.Range("D6").Value = 369
$ws.Range("D7").Value = 363
$ws.Range("D8").Value = 365
$ws.Range("D9").Value = 365
$ws.Range("D10").Value = 391
$ws.Range("D11").Value = 376
$ws.Range("D12").Value = 383
$ws.Range("D13").Value = 396
